# Applies the "records (95-04) import" edit to the usages_cp_format sheet.
#
# Summary of the change (derived from the OOXML diff):
#  - The "A"-section (legacy) rows had their min_year moved back from 2000
#    to 1995 (a few exceptions keep 2000/2011/2011 as before).
#  - Three brand-new usage rows are inserted/appended to the "A" section:
#      * "Refrigeration Manufacturing AC"  1995-1999  (row 7, inserted)
#      * "Fumigation"                      1995-2004  (inserted after
#                                           "Process agent")
#      * "Sterilant"                       2000-2001  (appended at the end)
#  - "Methyl bromide QPS" / "Methyl bromide Non-QPS" move from 2000 to 1998.
#  - "Tobacco fluffing" min_year moves from 2000 to 1995 (max_year 2011 kept).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert "Refrigeration Manufacturing AC" (1995-1999, A) as new row 7,
#    pushing the existing rows 7+ down by one.
# ---------------------------------------------------------------------
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Refrigeration Manufacturing AC"
$ws.Range("B7").Value = 1995
$ws.Range("C7").Value = 1999
$ws.Range("D7").Value = "A"

# ---------------------------------------------------------------------
# 2) min_year 2000 -> 1995 for the first five original "A" rows
#    (Aerosol, Foam, Fire fighting, Refrigeration, Refrigeration
#    Manufacturing), now sitting in rows 2-6.
# ---------------------------------------------------------------------
$ws.Range("B2:B6").Value = 1995

# ---------------------------------------------------------------------
# 3) min_year 2000 -> 1995 for "Refrigeration Servicing" (A) and
#    "Solvent application" (A, trailing tab) which are now rows 17-18.
# ---------------------------------------------------------------------
$ws.Range("B17").Value = 1995
$ws.Range("B18").Value = 1995

# ---------------------------------------------------------------------
# 4) min_year 2000 -> 1995 for "Process agent" (A), now row 22.
# ---------------------------------------------------------------------
$ws.Range("B22").Value = 1995

# ---------------------------------------------------------------------
# 5) Insert "Fumigation" (1995-2004, A) as new row 23, pushing MDI etc.
#    down by one.
# ---------------------------------------------------------------------
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "Fumigation"
$ws.Range("B23").Value = 1995
$ws.Range("C23").Value = 2004
$ws.Range("D23").Value = "A"

# ---------------------------------------------------------------------
# 6) min_year 2000 -> 1998 for "Methyl bromide QPS" / "Methyl bromide
#    Non-QPS", now rows 27-28, and min_year 2000 -> 1995 for "Tobacco
#    fluffing" (max_year 2011 unchanged), now row 29.
# ---------------------------------------------------------------------
$ws.Range("B27").Value = 1998
$ws.Range("B28").Value = 1998
$ws.Range("B29").Value = 1995

# ---------------------------------------------------------------------
# 7) Append "Sterilant" (2000-2001, A) as new row 30.
# ---------------------------------------------------------------------
$ws.Range("A30").Value = "Sterilant"
$ws.Range("B30").Value = 2000
$ws.Range("C30").Value = 2001
$ws.Range("D30").Value = "A"
